$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 1922.2222
$ws.Cells.Item(2, 9).Value = 1957
$ws.Cells.Item(2, 11).Value = 1957
$ws.Cells.Item(2, 13).Value = -1844
$ws.Cells.Item(28, 8).Value = 575
$ws.Cells.Item(28, 9).Value = 70.666664
$ws.Cells.Item(28, 10).Value = 877.6
$ws.Cells.Item(28, 11).Value = 70.666664
$ws.Cells.Item(28, 12).Value = 877.6
$ws.Cells.Item(28, 13).Value = 414.333336
$ws.Cells.Item(28, 14).Value = -1847.6
$ws.Cells.Item(33, 8).Value = 608
$ws.Cells.Item(33, 9).Value = 608
$ws.Cells.Item(33, 11).Value = 608
$ws.Cells.Item(33, 13).Value = -379
$ws.Cells.Item(42, 8).Value = 197.58824
$ws.Cells.Item(42, 9).Value = 136.2
$ws.Cells.Item(42, 10).Value = 285.2857
$ws.Cells.Item(42, 11).Value = 408.6
$ws.Cells.Item(42, 12).Value = 855.8571000000001
$ws.Cells.Item(42, 13).Value = -178.6
$ws.Cells.Item(42, 14).Value = -1315.8571
$ws.Cells.Item(43, 8).Value = 7499.5
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 7499.5
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 7499.5
$ws.Cells.Item(43, 13).ClearContents()
$ws.Cells.Item(43, 14).Value = -7637.5
$ws.Cells.Item(58, 8).Value = 91.666664
$ws.Cells.Item(58, 9).Value = 91.666664
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 274.999992
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = -124.999992
$ws.Cells.Item(58, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 404.57144
$ws.Cells.Item(61, 9).Value = 387.83334
$ws.Cells.Item(61, 10).Value = 505
$ws.Cells.Item(61, 11).Value = 1163.50002
$ws.Cells.Item(61, 12).Value = 1515
$ws.Cells.Item(61, 13).Value = -991.5000199999999
$ws.Cells.Item(61, 14).Value = -1859
$ws.Cells.Item(62, 8).Value = 6100
$ws.Cells.Item(62, 9).Value = 5600
$ws.Cells.Item(62, 10).Value = 7100
$ws.Cells.Item(62, 11).Value = 5600
$ws.Cells.Item(62, 12).Value = 7100
$ws.Cells.Item(62, 13).Value = -4976
$ws.Cells.Item(62, 14).Value = -8348
$ws.Cells.Item(65, 8).Value = 6100
$ws.Cells.Item(65, 9).Value = 5600
$ws.Cells.Item(65, 10).Value = 7100
$ws.Cells.Item(65, 11).Value = 28000
$ws.Cells.Item(65, 12).Value = 35500
$ws.Cells.Item(65, 13).Value = -24880
$ws.Cells.Item(65, 14).Value = -41740
$ws.Cells.Item(69, 8).Value = 30310770
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 30310770
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 90932310
$ws.Cells.Item(69, 13).ClearContents()
$ws.Cells.Item(69, 14).Value = -90934058
$ws.Cells.Item(72, 8).Value = 30310770
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 30310770
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 12).Value = 272796930
$ws.Cells.Item(72, 13).ClearContents()
$ws.Cells.Item(72, 14).Value = -272805666
$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 13).ClearContents()
$ws.Cells.Item(76, 14).ClearContents()
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 13).ClearContents()
$ws.Cells.Item(79, 14).ClearContents()
$ws.Cells.Item(86, 8).Value = 3950
$ws.Cells.Item(86, 9).Value = 4000
$ws.Cells.Item(86, 10).Value = 3900
$ws.Cells.Item(86, 11).Value = 4000
$ws.Cells.Item(86, 12).Value = 3900
$ws.Cells.Item(86, 13).Value = -2877
$ws.Cells.Item(86, 14).Value = -6146
$ws.Cells.Item(89, 8).Value = 3950
$ws.Cells.Item(89, 9).Value = 4000
$ws.Cells.Item(89, 10).Value = 3900
$ws.Cells.Item(89, 11).Value = 20000
$ws.Cells.Item(89, 12).Value = 19500
$ws.Cells.Item(89, 13).Value = -14384
$ws.Cells.Item(89, 14).Value = -30732
$ws.Cells.Item(103, 8).Value = 503.75
$ws.Cells.Item(103, 9).Value = 707.5
$ws.Cells.Item(103, 11).Value = 2122.5
$ws.Cells.Item(103, 13).Value = -1536.5
$ws.Cells.Item(106, 8).Value = 3804.25
$ws.Cells.Item(106, 9).Value = 3804.25
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 11).Value = 3804.25
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 13).Value = -3173.25
$ws.Cells.Item(106, 14).ClearContents()
$ws.Cells.Item(115, 8).Value = 450
$ws.Cells.Item(115, 9).Value = 450
$ws.Cells.Item(115, 11).Value = 1350
$ws.Cells.Item(115, 13).Value = 217
$ws.Cells.Item(135, 8).Value = 1003.5714
$ws.Cells.Item(135, 9).Value = 789.6667
$ws.Cells.Item(135, 10).Value = 1388.6
$ws.Cells.Item(135, 11).Value = 7107.0003
$ws.Cells.Item(135, 12).Value = 12497.4
$ws.Cells.Item(135, 13).Value = -4572.0003
$ws.Cells.Item(135, 14).Value = -17567.4
$ws.Cells.Item(137, 8).Value = 76564.21000000001
$ws.Cells.Item(137, 9).Value = 112772.81
$ws.Cells.Item(137, 10).Value = 4147
$ws.Cells.Item(137, 11).Value = 338318.43
$ws.Cells.Item(137, 12).Value = 12441
$ws.Cells.Item(137, 13).Value = -335768.43
$ws.Cells.Item(137, 14).Value = -17541
$ws.Cells.Item(138, 8).Value = 3037.6345
$ws.Cells.Item(138, 9).Value = 2096.625
$ws.Cells.Item(138, 10).Value = 3455.861
$ws.Cells.Item(138, 11).Value = 6289.875
$ws.Cells.Item(138, 12).Value = 10367.583
$ws.Cells.Item(138, 13).Value = -1149.875
$ws.Cells.Item(138, 14).Value = -20647.583

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 240.5
$ws.Cells.Item(4, 9).Value = 258.6
$ws.Cells.Item(4, 11).Value = 258.6
$ws.Cells.Item(4, 13).Value = -142.6
$ws.Cells.Item(32, 8).Value = 2704.1265
$ws.Cells.Item(32, 10).Value = 15144.333
$ws.Cells.Item(32, 12).Value = 15144.333
$ws.Cells.Item(32, 14).Value = -15718.333
$ws.Cells.Item(74, 8).Value = 56512.973
$ws.Cells.Item(74, 9).Value = 6893.387
$ws.Cells.Item(74, 10).Value = 364154.4
$ws.Cells.Item(74, 11).Value = 6893.387
$ws.Cells.Item(74, 12).Value = 364154.4
$ws.Cells.Item(74, 13).Value = -6019.387
$ws.Cells.Item(74, 14).Value = -365902.4
$ws.Cells.Item(77, 8).Value = 56512.973
$ws.Cells.Item(77, 9).Value = 6893.387
$ws.Cells.Item(77, 10).Value = 364154.4
$ws.Cells.Item(77, 11).Value = 34466.935
$ws.Cells.Item(77, 12).Value = 1820772
$ws.Cells.Item(77, 13).Value = -30098.935
$ws.Cells.Item(77, 14).Value = -1829508
$ws.Cells.Item(97, 8).Value = 9349.75
$ws.Cells.Item(97, 9).Value = 9349.75
$ws.Cells.Item(97, 11).Value = 9349.75
$ws.Cells.Item(97, 13).Value = -8853.75
$ws.Cells.Item(102, 8).Value = 5852.7144
$ws.Cells.Item(102, 9).Value = 4800
$ws.Cells.Item(102, 10).Value = 8484.5
$ws.Cells.Item(102, 11).Value = 4800
$ws.Cells.Item(102, 12).Value = 8484.5
$ws.Cells.Item(102, 13).Value = -3178
$ws.Cells.Item(102, 14).Value = -11728.5
$ws.Cells.Item(122, 8).Value = 2824.7144
$ws.Cells.Item(122, 9).Value = 2088.125
$ws.Cells.Item(122, 11).Value = 6264.375
$ws.Cells.Item(122, 13).Value = -3814.375
$ws.Cells.Item(132, 8).Value = 2344.9412
$ws.Cells.Item(132, 9).Value = 1759.1538
$ws.Cells.Item(132, 10).Value = 4248.75
$ws.Cells.Item(132, 11).Value = 5277.4614
$ws.Cells.Item(132, 12).Value = 12746.25
$ws.Cells.Item(132, 13).Value = -2747.4614
$ws.Cells.Item(132, 14).Value = -17806.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1542.3182
$ws.Cells.Item(20, 9).Value = 1346
$ws.Cells.Item(20, 11).Value = 1346
$ws.Cells.Item(20, 13).Value = -1099
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).ClearContents()
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(86, 8).Value = 4223.359
$ws.Cells.Item(86, 10).Value = 1612.8125
$ws.Cells.Item(86, 12).Value = 1612.8125
$ws.Cells.Item(86, 14).Value = -3858.8125
$ws.Cells.Item(89, 8).Value = 4223.359
$ws.Cells.Item(89, 10).Value = 1612.8125
$ws.Cells.Item(89, 12).Value = 8064.0625
$ws.Cells.Item(89, 14).Value = -19296.0625
$ws.Cells.Item(134, 8).Value = 6869.2666
$ws.Cells.Item(134, 9).Value = 3332.8333
$ws.Cells.Item(134, 10).Value = 9226.888999999999
$ws.Cells.Item(134, 11).Value = 9998.499899999999
$ws.Cells.Item(134, 12).Value = 27680.667
$ws.Cells.Item(134, 13).Value = -7463.499899999999
$ws.Cells.Item(134, 14).Value = -32750.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 266.5
$ws.Cells.Item(7, 9).Value = 38.07143
$ws.Cells.Item(7, 11).Value = 38.07143
$ws.Cells.Item(7, 13).Value = 74.92857000000001
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).ClearContents()
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(31, 8).Value = 93018.39999999999
$ws.Cells.Item(31, 9).Value = 2346
$ws.Cells.Item(31, 10).Value = 153466.67
$ws.Cells.Item(31, 11).Value = 2346
$ws.Cells.Item(31, 12).Value = 153466.67
$ws.Cells.Item(31, 13).Value = -2051
$ws.Cells.Item(31, 14).Value = -154056.67
$ws.Cells.Item(34, 8).Value = 93018.39999999999
$ws.Cells.Item(34, 9).Value = 2346
$ws.Cells.Item(34, 10).Value = 153466.67
$ws.Cells.Item(34, 11).Value = 2346
$ws.Cells.Item(34, 12).Value = 153466.67
$ws.Cells.Item(34, 13).Value = -2144
$ws.Cells.Item(34, 14).Value = -153870.67
$ws.Cells.Item(102, 8).Value = 44731
$ws.Cells.Item(102, 10).Value = 44731
$ws.Cells.Item(102, 12).Value = 44731
$ws.Cells.Item(102, 14).Value = -49599
$ws.Cells.Item(105, 8).Value = 1608.2084
$ws.Cells.Item(105, 9).Value = 1701.1428
$ws.Cells.Item(105, 10).Value = 957.6667
$ws.Cells.Item(105, 11).Value = 1701.1428
$ws.Cells.Item(105, 12).Value = 957.6667
$ws.Cells.Item(105, 13).Value = 45.85719999999992
$ws.Cells.Item(105, 14).Value = -4451.6667
$ws.Cells.Item(122, 8).Value = 3414.5
$ws.Cells.Item(122, 9).Value = 2970.6365
$ws.Cells.Item(122, 11).Value = 8911.9095
$ws.Cells.Item(122, 13).Value = -6461.9095
$ws.Cells.Item(132, 8).Value = 113516.625
$ws.Cells.Item(132, 10).Value = 299611
$ws.Cells.Item(132, 12).Value = 898833
$ws.Cells.Item(132, 14).Value = -903893
$ws.Cells.Item(134, 8).Value = 3891.6843
$ws.Cells.Item(134, 9).Value = 3228.4
$ws.Cells.Item(134, 10).Value = 4628.6665
$ws.Cells.Item(134, 11).Value = 9685.200000000001
$ws.Cells.Item(134, 12).Value = 13885.9995
$ws.Cells.Item(134, 13).Value = -7150.200000000001
$ws.Cells.Item(134, 14).Value = -18955.9995
$ws.Cells.Item(139, 8).Value = 140000
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 13).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 245.48979
$ws.Cells.Item(2, 9).Value = 122.833336
$ws.Cells.Item(2, 11).Value = 737.000016
$ws.Cells.Item(2, 13).Value = -624.000016
$ws.Cells.Item(68, 8).Value = 706.2
$ws.Cells.Item(68, 9).Value = 626
$ws.Cells.Item(68, 10).Value = 740.5714
$ws.Cells.Item(68, 11).Value = 1878
$ws.Cells.Item(68, 12).Value = 2221.7142
$ws.Cells.Item(68, 13).Value = -1067
$ws.Cells.Item(68, 14).Value = -3843.7142
$ws.Cells.Item(71, 8).Value = 706.2
$ws.Cells.Item(71, 9).Value = 626
$ws.Cells.Item(71, 10).Value = 740.5714
$ws.Cells.Item(71, 11).Value = 5634
$ws.Cells.Item(71, 12).Value = 6665.1426
$ws.Cells.Item(71, 13).Value = -1578
$ws.Cells.Item(71, 14).Value = -14777.1426
$ws.Cells.Item(113, 8).Value = 2497.5757
$ws.Cells.Item(113, 9).Value = 3623.6
$ws.Cells.Item(113, 10).Value = 2008
$ws.Cells.Item(113, 11).Value = 10870.8
$ws.Cells.Item(113, 12).Value = 6024
$ws.Cells.Item(113, 13).Value = -8700.799999999999
$ws.Cells.Item(113, 14).Value = -10364
$ws.Cells.Item(122, 8).Value = 1000.4545
$ws.Cells.Item(122, 9).Value = 918.7143
$ws.Cells.Item(122, 10).Value = 1143.5
$ws.Cells.Item(122, 11).Value = 8268.4287
$ws.Cells.Item(122, 12).Value = 10291.5
$ws.Cells.Item(122, 13).Value = -5818.4287
$ws.Cells.Item(122, 14).Value = -15191.5
$ws.Cells.Item(132, 8).Value = 2079.8125
$ws.Cells.Item(132, 9).Value = 1275
$ws.Cells.Item(132, 10).Value = 2445.6365
$ws.Cells.Item(132, 11).Value = 11475
$ws.Cells.Item(132, 12).Value = 22010.7285
$ws.Cells.Item(132, 13).Value = -8945
$ws.Cells.Item(132, 14).Value = -27070.7285
$ws.Cells.Item(133, 8).Value = 2693.111
$ws.Cells.Item(133, 9).Value = 2405.375
$ws.Cells.Item(133, 11).Value = 7216.125
$ws.Cells.Item(133, 13).Value = -2156.125
$ws.Cells.Item(140, 8).Value = 2991.25
$ws.Cells.Item(140, 9).Value = 2991.25
$ws.Cells.Item(140, 11).Value = 8973.75
$ws.Cells.Item(140, 13).Value = -3793.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 4606.75
$ws.Cells.Item(80, 10).Value = 20006
$ws.Cells.Item(80, 12).Value = 20006
$ws.Cells.Item(80, 14).Value = -22002
$ws.Cells.Item(83, 8).Value = 4606.75
$ws.Cells.Item(83, 10).Value = 20006
$ws.Cells.Item(83, 12).Value = 100030
$ws.Cells.Item(83, 14).Value = -110014
$ws.Cells.Item(97, 8).Value = 1228
$ws.Cells.Item(97, 9).Value = 1228
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 1228
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = -732
$ws.Cells.Item(97, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 498171.56
$ws.Cells.Item(122, 9).Value = 743840.8
$ws.Cells.Item(122, 11).Value = 2231522.4
$ws.Cells.Item(122, 13).Value = -2229072.4
$ws.Cells.Item(132, 8).Value = 4551.8125
$ws.Cells.Item(132, 9).Value = 3569.2856
$ws.Cells.Item(132, 10).Value = 5316
$ws.Cells.Item(132, 11).Value = 10707.8568
$ws.Cells.Item(132, 12).Value = 15948
$ws.Cells.Item(132, 13).Value = -8177.856800000001
$ws.Cells.Item(132, 14).Value = -21008

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(47, 8).Value = 30000
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 13).ClearContents()
$ws.Cells.Item(52, 8).Value = 30000
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 13).ClearContents()
$ws.Cells.Item(55, 8).Value = 1916.2307
$ws.Cells.Item(55, 9).Value = 1526.0667
$ws.Cells.Item(55, 11).Value = 1526.0667
$ws.Cells.Item(55, 13).Value = -1353.0667
$ws.Cells.Item(93, 8).Value = 4180.722
$ws.Cells.Item(93, 9).Value = 3088.1428
$ws.Cells.Item(93, 10).Value = 8004.75
$ws.Cells.Item(93, 11).Value = 3088.1428
$ws.Cells.Item(93, 12).Value = 8004.75
$ws.Cells.Item(93, 13).Value = -1840.1428
$ws.Cells.Item(93, 14).Value = -10500.75
$ws.Cells.Item(100, 8).Value = 3250.1333
$ws.Cells.Item(100, 9).Value = 3125.1428
$ws.Cells.Item(100, 11).Value = 3125.1428
$ws.Cells.Item(100, 13).Value = -2584.1428
$ws.Cells.Item(106, 8).Value = 9444
$ws.Cells.Item(106, 10).Value = 9444
$ws.Cells.Item(106, 12).Value = 9444
$ws.Cells.Item(106, 14).Value = -11968
$ws.Cells.Item(122, 8).Value = 7368.1
$ws.Cells.Item(122, 9).Value = 4475.75
$ws.Cells.Item(122, 11).Value = 13427.25
$ws.Cells.Item(122, 13).Value = -10977.25
$ws.Cells.Item(132, 8).Value = 4680.55
$ws.Cells.Item(132, 9).Value = 3135.6667
$ws.Cells.Item(132, 10).Value = 5342.643
$ws.Cells.Item(132, 11).Value = 9407.000100000001
$ws.Cells.Item(132, 12).Value = 16027.929
$ws.Cells.Item(132, 13).Value = -6877.000100000001
$ws.Cells.Item(132, 14).Value = -21087.929
$ws.Cells.Item(136, 8).Value = 66130.44
$ws.Cells.Item(136, 9).Value = 70272.13
$ws.Cells.Item(136, 10).Value = 4005
$ws.Cells.Item(136, 11).Value = 210816.39
$ws.Cells.Item(136, 12).Value = 12015
$ws.Cells.Item(136, 13).Value = -208266.39
$ws.Cells.Item(136, 14).Value = -17115

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(37, 8).Value = 35562
$ws.Cells.Item(37, 9).Value = 35000
$ws.Cells.Item(37, 11).Value = 35000
$ws.Cells.Item(37, 13).Value = -34797
$ws.Cells.Item(62, 8).Value = 8608.200000000001
$ws.Cells.Item(62, 9).Value = 2800
$ws.Cells.Item(62, 10).Value = 8850.208000000001
$ws.Cells.Item(62, 11).Value = 2800
$ws.Cells.Item(62, 12).Value = 8850.208000000001
$ws.Cells.Item(62, 13).Value = -2176
$ws.Cells.Item(62, 14).Value = -10098.208
$ws.Cells.Item(65, 8).Value = 8608.200000000001
$ws.Cells.Item(65, 9).Value = 2800
$ws.Cells.Item(65, 10).Value = 8850.208000000001
$ws.Cells.Item(65, 11).Value = 14000
$ws.Cells.Item(65, 12).Value = 44251.04
$ws.Cells.Item(65, 13).Value = -10880
$ws.Cells.Item(65, 14).Value = -50491.04
$ws.Cells.Item(107, 8).Value = 3966
$ws.Cells.Item(107, 9).Value = 4714.091
$ws.Cells.Item(107, 10).Value = 1223
$ws.Cells.Item(107, 11).Value = 14142.273
$ws.Cells.Item(107, 12).Value = 3669
$ws.Cells.Item(107, 13).Value = -12222.273
$ws.Cells.Item(107, 14).Value = -7509
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 749.6070999999999
$ws.Cells.Item(113, 9).Value = 550.7895
$ws.Cells.Item(113, 10).Value = 1169.3334
$ws.Cells.Item(113, 11).Value = 1652.3685
$ws.Cells.Item(113, 12).Value = 3508.0002
$ws.Cells.Item(113, 13).Value = 517.6315
$ws.Cells.Item(113, 14).Value = -7848.0002
$ws.Cells.Item(114, 8).Value = 49991.5
$ws.Cells.Item(114, 10).Value = 49991.5
$ws.Cells.Item(114, 12).Value = 49991.5
$ws.Cells.Item(114, 14).Value = -58669.5
$ws.Cells.Item(116, 8).Value = 40999.5
$ws.Cells.Item(116, 10).Value = 40999.5
$ws.Cells.Item(116, 12).Value = 40999.5
$ws.Cells.Item(116, 14).Value = -50177.5
$ws.Cells.Item(117, 8).Value = 22000
$ws.Cells.Item(117, 10).Value = 22000
$ws.Cells.Item(117, 12).Value = 22000
$ws.Cells.Item(117, 14).Value = -31178
$ws.Cells.Item(119, 8).Value = 41099.668
$ws.Cells.Item(119, 10).Value = 41099.668
$ws.Cells.Item(119, 12).Value = 41099.668
$ws.Cells.Item(119, 14).Value = -50775.668
$ws.Cells.Item(121, 8).Value = 54999
$ws.Cells.Item(121, 10).Value = 54999
$ws.Cells.Item(121, 12).Value = 54999
$ws.Cells.Item(121, 14).Value = -58493
$ws.Cells.Item(122, 8).Value = 4995
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 4995
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).Value = 14985
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(122, 14).Value = -19885
$ws.Cells.Item(126, 8).Value = 2119.7273
$ws.Cells.Item(126, 9).Value = 2211.2
$ws.Cells.Item(126, 11).Value = 6633.599999999999
$ws.Cells.Item(126, 13).Value = -4163.599999999999
$ws.Cells.Item(132, 8).Value = 131886.42
$ws.Cells.Item(132, 10).Value = 302279
$ws.Cells.Item(132, 12).Value = 906837
$ws.Cells.Item(132, 14).Value = -911897
$ws.Cells.Item(135, 8).Value = 72622.875
$ws.Cells.Item(135, 9).Value = 55499.832
$ws.Cells.Item(135, 10).Value = 123992
$ws.Cells.Item(135, 11).Value = 55499.832
$ws.Cells.Item(135, 12).Value = 123992
$ws.Cells.Item(135, 13).Value = -50429.832
$ws.Cells.Item(135, 14).Value = -134132
$ws.Cells.Item(136, 8).Value = 5027.3335
$ws.Cells.Item(136, 9).Value = 3751
$ws.Cells.Item(136, 10).Value = 5665.5
$ws.Cells.Item(136, 11).Value = 11253
$ws.Cells.Item(136, 12).Value = 16996.5
$ws.Cells.Item(136, 13).Value = -8703
$ws.Cells.Item(136, 14).Value = -22096.5
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 14).ClearContents()
